# Scheduled marketboard refresh for Sheets/Typhon_Profits.xlsx
# Updates currentAveragePrice*/Leve cost & profit columns (H:N) on each crafting-class
# worksheet with freshly polled prices. A couple of rows lose their trailing (N)
# profit cell entirely where the refreshed recipe no longer has a second yield path.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 543
$ws.Range("H40").Value = 1128.0385
$ws.Range("I40").Value = 796.4286
$ws.Range("J40").Value = 1514.9166
$ws.Range("K40").Value = 796.4286
$ws.Range("L40").Value = 1514.9166
$ws.Range("M40").Value = -621.4286
$ws.Range("N40").Value = -1864.9166
$ws.Range("H64").Value = 3159.9
$ws.Range("I64").Value = 2924.75
$ws.Range("J64").Value = 3316.6667
$ws.Range("K64").Value = 2924.75
$ws.Range("L64").Value = 3316.6667
$ws.Range("M64").Value = -2676.75
$ws.Range("N64").Value = -3812.6667
$ws.Range("H67").Value = 3159.9
$ws.Range("I67").Value = 2924.75
$ws.Range("J67").Value = 3316.6667
$ws.Range("K67").Value = 2924.75
$ws.Range("L67").Value = 3316.6667
$ws.Range("M67").Value = -2066.75
$ws.Range("N67").Value = -5032.6667
$ws.Range("H70").Value = 1243.1666
$ws.Range("I70").Value = 1259.6666
$ws.Range("K70").Value = 3778.9998
$ws.Range("M70").Value = -3508.9998
$ws.Range("H73").Value = 1243.1666
$ws.Range("I73").Value = 1259.6666
$ws.Range("K73").Value = 3778.9998
$ws.Range("M73").Value = -2842.9998
$ws.Range("H107").Value = 1057.2106
$ws.Range("I107").Value = 812.5333000000001
$ws.Range("J107").Value = 1974.75
$ws.Range("K107").Value = 812.5333000000001
$ws.Range("L107").Value = 1974.75
$ws.Range("M107").Value = 1107.4667
$ws.Range("N107").Value = -5814.75
$ws.Range("H138").Value = 136009.77
$ws.Range("I138").Value = 3024.875
$ws.Range("J138").Value = 151655.05
$ws.Range("K138").Value = 9074.625
$ws.Range("L138").Value = 454965.15
$ws.Range("M138").Value = -3934.625
$ws.Range("N138").Value = -465245.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1877.2954
$ws.Range("I61").Value = 1550.1
$ws.Range("J61").Value = 5149.25
$ws.Range("K61").Value = 1550.1
$ws.Range("L61").Value = 5149.25
$ws.Range("M61").Value = -1338.1
$ws.Range("N61").Value = -5573.25
$ws.Range("H63").Value = 2406186.2
$ws.Range("I63").Value = 2601.818
$ws.Range("K63").Value = 2601.818
$ws.Range("M63").Value = -1915.818
$ws.Range("H66").Value = 2406186.2
$ws.Range("I66").Value = 2601.818
$ws.Range("K66").Value = 13009.09
$ws.Range("M66").Value = -9577.09
$ws.Range("H97").Value = 555.25
$ws.Range("I97").Value = 572.5
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 572.5
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = -76.5
$ws.Range("N97").Value = -1392
$ws.Range("H110").Value = 1192.4166
$ws.Range("I110").Value = 1166.2222
$ws.Range("K110").Value = 1166.2222
$ws.Range("M110").Value = 878.7778000000001
$ws.Range("H122").Value = 2252.68
$ws.Range("I122").Value = 2164.682
$ws.Range("J122").Value = 2898
$ws.Range("K122").Value = 6494.045999999999
$ws.Range("L122").Value = 8694
$ws.Range("M122").Value = -4044.045999999999
$ws.Range("N122").Value = -13594
$ws.Range("H132").Value = 12801.106
$ws.Range("I132").Value = 1956.5
$ws.Range("K132").Value = 5869.5
$ws.Range("M132").Value = -3339.5
$ws.Range("H136").Value = 1877.2954
$ws.Range("I136").Value = 1550.1
$ws.Range("J136").Value = 5149.25
$ws.Range("K136").Value = 4650.299999999999
$ws.Range("L136").Value = 15447.75
$ws.Range("M136").Value = -2100.299999999999
$ws.Range("N136").Value = -20547.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1988.3684
$ws.Range("I86").Value = 1853.7333
$ws.Range("J86").Value = 2493.25
$ws.Range("K86").Value = 1853.7333
$ws.Range("L86").Value = 2493.25
$ws.Range("M86").Value = -730.7333000000001
$ws.Range("N86").Value = -4739.25
$ws.Range("H89").Value = 1988.3684
$ws.Range("I89").Value = 1853.7333
$ws.Range("J89").Value = 2493.25
$ws.Range("K89").Value = 9268.666500000001
$ws.Range("L89").Value = 12466.25
$ws.Range("M89").Value = -3652.666500000001
$ws.Range("N89").Value = -23698.25
$ws.Range("H94").Value = 928.48
$ws.Range("I94").Value = 681.13336
$ws.Range("K94").Value = 681.13336
$ws.Range("M94").Value = -230.13336
$ws.Range("H97").Value = 13549.214
$ws.Range("I97").Value = 5269.857
$ws.Range("J97").Value = 21828.572
$ws.Range("K97").Value = 5269.857
$ws.Range("L97").Value = 21828.572
$ws.Range("M97").Value = -4278.857
$ws.Range("N97").Value = -23810.572
$ws.Range("H134").Value = 4178
$ws.Range("J134").Value = 3000
$ws.Range("L134").Value = 9000
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3290.4917
$ws.Range("I31").Value = 1521.2222
$ws.Range("J31").Value = 5838.24
$ws.Range("K31").Value = 1521.2222
$ws.Range("L31").Value = 5838.24
$ws.Range("M31").Value = -1226.2222
$ws.Range("N31").Value = -6428.24
$ws.Range("H34").Value = 3290.4917
$ws.Range("I34").Value = 1521.2222
$ws.Range("J34").Value = 5838.24
$ws.Range("K34").Value = 1521.2222
$ws.Range("L34").Value = 5838.24
$ws.Range("M34").Value = -1319.2222
$ws.Range("N34").Value = -6242.24
$ws.Range("H50").Value = 14285.714
$ws.Range("J50").Value = 14285.714
$ws.Range("L50").Value = 14285.714
$ws.Range("N50").Value = -15535.714
$ws.Range("H51").Value = 34666.668
$ws.Range("J51").Value = 34666.668
$ws.Range("L51").Value = 34666.668
$ws.Range("N51").Value = -36138.668
$ws.Range("H59").Value = 21150
$ws.Range("J59").Value = 26533.334
$ws.Range("L59").Value = 26533.334
$ws.Range("N59").Value = -28823.334
$ws.Range("H60").Value = 17870.428
$ws.Range("J60").Value = 21000
$ws.Range("L60").Value = 21000
$ws.Range("N60").Value = -22022
$ws.Range("H61").Value = 34666.668
$ws.Range("J61").Value = 34666.668
$ws.Range("L61").Value = 34666.668
$ws.Range("N61").Value = -35362.668
$ws.Range("H97").Value = 32197
$ws.Range("J97").Value = 32197
$ws.Range("L97").Value = 32197
$ws.Range("N97").Value = -34179
$ws.Range("H99").Value = 19234420
$ws.Range("I99").Value = 3182.9092
$ws.Range("J99").Value = 33337326
$ws.Range("K99").Value = 3182.9092
$ws.Range("L99").Value = 33337326
$ws.Range("M99").Value = -1684.9092
$ws.Range("N99").Value = -33340322
$ws.Range("H126").Value = 19234420
$ws.Range("I126").Value = 3182.9092
$ws.Range("J126").Value = 33337326
$ws.Range("K126").Value = 9548.7276
$ws.Range("L126").Value = 100011978
$ws.Range("M126").Value = -7078.7276
$ws.Range("N126").Value = -100016918
$ws.Range("H134").Value = 1013.9355
$ws.Range("I134").Value = 814.8929000000001
$ws.Range("K134").Value = 2444.6787
$ws.Range("M134").Value = 90.32129999999961

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 894.6667
$ws.Range("J113").Value = 943.3333
$ws.Range("L113").Value = 2829.9999
$ws.Range("N113").Value = -7169.9999
$ws.Range("H131").Value = 645.28
$ws.Range("J131").Value = 769.31506
$ws.Range("L131").Value = 2307.94518
$ws.Range("N131").Value = -12387.94518
$ws.Range("H140").Value = 2281.7896
$ws.Range("I140").Value = 1306.3572
$ws.Range("J140").Value = 5013
$ws.Range("K140").Value = 3919.0716
$ws.Range("L140").Value = 15039
$ws.Range("M140").Value = 1260.9284
$ws.Range("N140").Value = -25399
$ws.Range("H141").Value = 5740
$ws.Range("I141").Value = 5740
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 17220
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -12040
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2757.2
$ws.Range("I102").Value = 2478.2354
$ws.Range("J102").Value = 4338
$ws.Range("K102").Value = 2478.2354
$ws.Range("L102").Value = 4338
$ws.Range("M102").Value = -856.2354
$ws.Range("N102").Value = -7582
$ws.Range("H126").Value = 5378.2173
$ws.Range("J126").Value = 4745.364
$ws.Range("L126").Value = 14236.092
$ws.Range("N126").Value = -19176.092
$ws.Range("H132").Value = 14924.637
$ws.Range("I132").Value = 4291.909
$ws.Range("K132").Value = 12875.727
$ws.Range("M132").Value = -10345.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2266
$ws.Range("I68").Value = 1899.5
$ws.Range("K68").Value = 1899.5
$ws.Range("M68").Value = -1150.5
$ws.Range("H71").Value = 2266
$ws.Range("I71").Value = 1899.5
$ws.Range("K71").Value = 9497.5
$ws.Range("M71").Value = -5753.5
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 625.6222
$ws.Range("I132").Value = 504.84375
$ws.Range("K132").Value = 1514.53125
$ws.Range("M132").Value = 1015.46875
